# Generate Report for Handback
# Fills in the Correspond Handoff/Handback datetimes that are produced once
# the handback round-trip for the second file (4c55c43b-...) completes, and
# propagates the latest timestamp back up to the Overview sheet.

$wb = $excel.ActiveWorkbook

# --- zh-cn sheet: row 3 is the 4c55c43b-...md / ...zh-cn.xlf entry ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H3").Value = "2016-10-24 06:34:24"
$wsZhCn.Range("K3").Value = "2016-10-24 06:35:06"

# --- de-de sheet: row 3 is the 4c55c43b-...md / ...de-de.xlf entry ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H3").Value = "2016-10-24 06:34:35"
$wsDeDe.Range("K3").Value = "2016-10-24 06:35:22"

# --- Overview sheet: row 3 (4c55c43b-...md) now shows the newer of the
#     two "Latest HO Xliff Generate Date" timestamps (the de-de handoff). ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G3").Value = "2016-10-24 06:34:35"
